$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.584.23'
$ws.Range('E2').Value = '  +2.91%  '

$ws.Range('D3').Value = '1.606.58'
$ws.Range('E3').Value = '  +2.51%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('E5').Value = '  +1.17%  '

$ws.Range('E6').Value = '  +5.27%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('E8').Value = '  +7.01%  '

$ws.Range('E9').Value = '  -1.46%  '

$ws.Range('E10').Value = '  +2.19%  '

$ws.Range('E11').Value = '  +2.68%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0908'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.48%  '

$ws.Range('D13').Value = '1.835.00'
$ws.Range('E13').Value = '  +2.48%  '

$ws.Range('D14').Value = '1.602.70'
$ws.Range('E14').Value = '  +2.40%  '

$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '29.586.24'
$ws.Range('E15').Value = '  +2.98%  '

$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.538'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.81%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.72'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.96%  '

$ws.Range('E18').Value = '  +3.29%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.41%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.36%  '

$ws.Range('E21').Value = '  +1.74%  '

$ws.Range('E22').Value = '  +0.10%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.99'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.50%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.23'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.72%  '

$ws.Range('E25').Value = '  +0.46%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.42'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.94%  '

$ws.Range('E27').Value = '  +3.93%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.91%  '

$ws.Range('E29').Value = '  +2.41%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.07%  '

$ws.Range('E31').Value = '  +2.53%  '

$ws.Range('E32').Value = '  +1.01%  '

$ws.Range('E33').Value = '  +1.28%  '

$ws.Range('E34').Value = '  +3.28%  '

$ws.Range('D35').Value = '1.409.56'
$ws.Range('E35').Value = '  +0.83%  '

$ws.Range('E36').Value = '  +0.47%  '

$ws.Range('E37').Value = '  +4.21%  '

$ws.Range('E38').Value = '  +3.42%  '

$ws.Range('E39').Value = '  +0.26%  '

$ws.Range('E40').Value = '  +2.29%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.539'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.68%  '

$ws.Range('E42').Value = '  +0.83%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0492'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.52%  '

$ws.Range('E44').Value = '  +24.34%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.799'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.26%  '

$ws.Range('E46').Value = '  +0.06%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.88%  '

$ws.Range('E48').Value = '  +0.84%  '

$ws.Range('D49').Value = '1.746.07'
$ws.Range('E49').Value = '  +2.71%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.855'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.82%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '86.75'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.91%  '
